$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Automatable Tests")

# check border info on existing cell D86 (style 39, borderId=5)
$c = $ws.Cells.Item(86,4)
try {
  Write-Output ("EdgeBottom LineStyle: " + $c.Borders.Item(9).LineStyle)
} catch { Write-Output ("ERR: " + $_.Exception.Message) }
try {
  Write-Output ("Borders(xlEdgeLeft).Weight: " + $c.Borders.Item(7).Weight)
} catch { Write-Output ("ERR2: " + $_.Exception.Message) }

# Try Range.Copy then regular Paste (values+formats) into row 87 after insert
$ws.Rows.Item(87).Insert()
$src = $ws.Range("A86:M86")
$dst = $ws.Range("A87:M87")
$src.Copy($dst)
Write-Output "Copy(dest) worked"
